$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '28.384.04'
$ws.Range("E2").Value = '  +3.36%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '1.867.95'
$ws.Range("E3").Value = '  +1.67%  '

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.14%  '

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4690'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.75%  '

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3957'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.71%  '

# Row 9 - OKB
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.17'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.08%  '

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08005'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.80%  '

# Row 11 - Polygon
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9983'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.44%  '

# Row 12 - Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.82'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.30%  '

# Row 13 - was WrappedEther, now Polkadot
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.991'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.65%  '

# Row 14 - was Polkadot, now WrappedEther
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.859.75'
$ws.Range("E14").Value = '  +1.18%  '

# Row 15 - Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.223'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.87%  '

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.73%  '

# Row 17 - BinanceUSD
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.11%  '

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001037'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.71%  '

# Row 19 - TRON
$ws.Range("E19").Value = '  -0.39%  '

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.36%  '

# Row 21 - Dai
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.01%  '

# Row 22 - WrappedBTC
$ws.Range("D22").Value = '28.393.09'
$ws.Range("E22").Value = '  +3.39%  '

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.446'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.07%  '

# Row 24 - Cosmos
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.04'
$ws.Range("D24").Style = "Normal"

# Row 25 - Toncoin
$ws.Range("E25").Value = '  -1.20%  '

# Row 26 - WrappedliquidstakedEther2.0
$ws.Range("D26").Value = '2.082.00'
$ws.Range("E26").Value = '  +1.15%  '

# Row 27 - Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.49'
$ws.Range("D27").Style = "Normal"

# Row 28 - EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.92%  '

# Row 29 - LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.125'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.88%  '

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.483'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.94%  '

# Row 31 - BitcoinCash
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '120.10'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.25%  '

# Row 32 - ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9649'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.19%  '

# Row 33 - Stellar
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09474'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.96%  '

# Row 34 - HuobiToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.574'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.12%  '

# Row 35 - Filecoin
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.340'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.21%  '

# Row 36 - ARBITRUM
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.370'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.83%  '

# Row 37 - Hedera
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06081'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.48%  '

# Row 38 - VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02241'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.75%  '

# Row 39 - FraxShare
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.360'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.59%  '

# Row 40 - TrustWalletToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.187'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.58%  '

# Row 41 - TheSandbox
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5929'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.10%  '

# Row 42 - Frax
$ws.Range("E42").Value = '  -0.03%  '

# Row 43 - Algorand
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1866'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.51%  '

# Row 44 - Aptos
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.80%  '

# Row 45 - WEMIXTOKEN
$ws.Range("E45").Value = '  +4.24%  '

# Row 46 - Decentraland
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5571'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.58%  '

# Row 47 - EnergySwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.02%  '

# Row 48 - NEARProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.952'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.64%  '

# Row 49 - Cronos
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06859'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.11%  '

# Row 50 - RenderToken
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.044'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +15.65%  '

# Row 51 - Quant
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '111.30'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.40%  '
